$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on all cells being updated so that numeric-looking
# strings (e.g. "1.000", "5.960") are preserved exactly as text and are not
# silently re-interpreted by Excel as numbers (which would drop trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.476.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.913.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4784'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4098'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.69'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08034'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.42'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.883.66'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.960'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.174'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.48'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.75'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.495.97'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.204'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.121.65'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.790'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.07%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.64'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.066'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09565'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.424'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.395'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06103'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02256'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.358'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.177'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1842'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.413'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07763'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.43%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5556'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.13'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.45'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.84'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.39%  '
